# This workbook's data rows (2-10) were reshuffled into a new order, and the
# "Ost"/"Nord" (Q/R) coordinate columns were rounded from precise computed
# floating point values down to plain integers.
#
# Row permutation (new row <- old row), derived by matching the "Id" column:
#   2 <- 8    3 <- 6    4 <- 10   5 <- 3    6 <- 2
#   7 <- 9    8 <- 7    9 <- 4    10 <- 5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns I (Antal), Y/AA (Start-/Slutdatum) hold values that look numeric or
# date-like ("38", "2023-09-16", ...) but must stay plain text, matching the
# source file. Force Text formatting before writing so Excel does not
# auto-convert them into numbers/dates.
$ws.Range("I2:I10").NumberFormat = "@"
$ws.Range("Y2:Y10").NumberFormat = "@"
$ws.Range("AA2:AA10").NumberFormat = "@"

# Snapshot every data row (whole used column range) before overwriting anything.
$row2 = $ws.Range("A2:AY2").Value()
$row3 = $ws.Range("A3:AY3").Value()
$row4 = $ws.Range("A4:AY4").Value()
$row5 = $ws.Range("A5:AY5").Value()
$row6 = $ws.Range("A6:AY6").Value()
$row7 = $ws.Range("A7:AY7").Value()
$row8 = $ws.Range("A8:AY8").Value()
$row9 = $ws.Range("A9:AY9").Value()
$row10 = $ws.Range("A10:AY10").Value()

# Write rows back out in their new order.
$ws.Range("A2:AY2").Value = $row8
$ws.Range("A3:AY3").Value = $row6
$ws.Range("A4:AY4").Value = $row10
$ws.Range("A5:AY5").Value = $row3
$ws.Range("A6:AY6").Value = $row2
$ws.Range("A7:AY7").Value = $row9
$ws.Range("A8:AY8").Value = $row7
$ws.Range("A9:AY9").Value = $row4
$ws.Range("A10:AY10").Value = $row5

# Round the Ost (Q) / Nord (R) coordinates to whole numbers for every data row.
for ($r = 2; $r -le 10; $r++) {
    $ost = $ws.Cells.Item($r, 17).Value()
    $ws.Cells.Item($r, 17).Value = [Math]::Round([double]$ost)

    $nord = $ws.Cells.Item($r, 18).Value()
    $ws.Cells.Item($r, 18).Value = [Math]::Round([double]$nord)
}
